# Test Suite punctuation cases.
# Insert two new checklist rows ("Same word" / "Same word, different case"
# under Punctuation > Case-insensitive), each annotated with the comment
# "Duplicate longest word returned only once.", and mark the three
# Case-insensitive / Spaces / Periods checks as PASS instead of Not Tested.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 15, shifting the remaining checklist
# rows (and the trailing formatted spacer row) down by two.
$ws.Rows("15:16").Insert()

# New sub-item rows under "Punctuation" > "Case-insensitive".
$ws.Range("D15").Value = "Same word"
$ws.Range("K15").Value = "Not Tested"
$ws.Range("M15").Value = "Duplicate longest word returned only once."

$ws.Range("D16").Value = "Same word, different case"
$ws.Range("K16").Value = "Not Tested"
$ws.Range("M16").Value = "Duplicate longest word returned only once."

# Mark these checklist items as passing.
$ws.Range("K4").Value = "PASS"
$ws.Range("K5").Value = "PASS"
$ws.Range("K7").Value = "PASS"
